$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 1465
$ws.Range("E2").Value = -12
$ws.Range("F2").Value = -12
$ws.Range("G2").Value = -24
$ws.Range("H2").Value = -27
$ws.Range("I2").Value = -27
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 1569
$ws.Range("L2").Value = 880
$ws.Range("M2").Value = 688
$ws.Range("N2").Value = 688
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 386
$ws.Range("Q2").Value = -147
$ws.Range("R2").Value = -144
$ws.Range("S2").Value = 271
$ws.Range("T2").Value = 145
$ws.Range("U2").Value = -292
$ws.Range("V2").Value = 553
$ws.Range("W2").Value = -0.83
$ws.Range("X2").Value = -1.84
$ws.Range("Y2").Value = -3.84
$ws.Range("Z2").Value = -1.91
$ws.Range("AA2").Value = 127.86
$ws.Range("AB2").Value = 76.20999999999999
$ws.Range("AC2").Value = -35
$ws.Range("AD2").Value = -22.94
$ws.Range("AE2").Value = 892
$ws.Range("AF2").Value = 0.9
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 77124820

# Row 3
$ws.Range("D3").Value = 1547
$ws.Range("E3").Value = -75
$ws.Range("F3").Value = -81
$ws.Range("G3").Value = 25
$ws.Range("H3").Value = 31
$ws.Range("I3").Value = 31
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 1672
$ws.Range("L3").Value = 986
$ws.Range("M3").Value = 686
$ws.Range("N3").Value = 686
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = 386
$ws.Range("Q3").Value = -110
$ws.Range("R3").Value = -31
$ws.Range("S3").Value = 158
$ws.Range("T3").Value = 56
$ws.Range("U3").Value = -166
$ws.Range("V3").Value = 745
$ws.Range("W3").Value = -4.86
$ws.Range("X3").Value = 2
$ws.Range("Y3").Value = 4.47
$ws.Range("Z3").Value = 1.91
$ws.Range("AA3").Value = 143.71
$ws.Range("AB3").Value = 83.31999999999999
$ws.Range("AC3").Value = 40
$ws.Range("AD3").Value = 36.68
$ws.Range("AE3").Value = 889
$ws.Range("AF3").Value = 1.64
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 77124820

# Row 4
$ws.Range("D4").Value = 1676
$ws.Range("E4").Value = 49
$ws.Range("F4").Value = 49
$ws.Range("G4").Value = 56
$ws.Range("H4").Value = 28
$ws.Range("I4").Value = 28
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 1601
$ws.Range("L4").Value = 894
$ws.Range("M4").Value = 706
$ws.Range("N4").Value = 706
$ws.Range("P4").Value = 386
$ws.Range("Q4").Value = 67
$ws.Range("R4").Value = 72
$ws.Range("S4").Value = -117
$ws.Range("T4").Value = 28
$ws.Range("U4").Value = 39
$ws.Range("V4").Value = 667
$ws.Range("W4").Value = 2.91
$ws.Range("X4").Value = 1.69
$ws.Range("Y4").Value = 4.08
$ws.Range("Z4").Value = 1.73
$ws.Range("AA4").Value = 126.62
$ws.Range("AB4").Value = 90.58
$ws.Range("AC4").Value = 37
$ws.Range("AD4").Value = 31.14
$ws.Range("AE4").Value = 916
$ws.Range("AF4").Value = 1.25
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 77124820
$ws.Range("O4").ClearContents()
$ws.Range("AG4").ClearContents()
$ws.Range("AH4").ClearContents()

# Row 5
$ws.Range("D5").Value = 1763
$ws.Range("E5").Value = 35
$ws.Range("F5").Value = 35
$ws.Range("G5").Value = -9
$ws.Range("H5").Value = -17
$ws.Range("I5").Value = -17
$ws.Range("K5").Value = 1398
$ws.Range("L5").Value = 705
$ws.Range("M5").Value = 693
$ws.Range("N5").Value = 693
$ws.Range("P5").Value = 386
$ws.Range("Q5").Value = 21
$ws.Range("R5").Value = 102
$ws.Range("S5").Value = -138
$ws.Range("T5").Value = 57
$ws.Range("U5").Value = -36
$ws.Range("V5").Value = 488
$ws.Range("W5").Value = 1.99
$ws.Range("X5").Value = -0.98
$ws.Range("Y5").Value = -2.47
$ws.Range("Z5").Value = -1.15
$ws.Range("AA5").Value = 101.76
$ws.Range("AB5").Value = 85.63
$ws.Range("AC5").Value = -22
$ws.Range("AD5").Value = -98.13
$ws.Range("AE5").Value = 898
$ws.Range("AF5").Value = 2.45
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 77124820
$ws.Range("J5").ClearContents()
$ws.Range("O5").ClearContents()
$ws.Range("AG5").ClearContents()
$ws.Range("AH5").ClearContents()

# Row 6
$ws.Range("D6").Value = 1923
$ws.Range("E6").Value = -120
$ws.Range("F6").Value = -120
$ws.Range("G6").Value = -135
$ws.Range("H6").Value = -139
$ws.Range("I6").Value = -139
$ws.Range("K6").Value = 1395
$ws.Range("L6").Value = 845
$ws.Range("M6").Value = 550
$ws.Range("N6").Value = 550
$ws.Range("P6").Value = 386
$ws.Range("Q6").Value = -43
$ws.Range("R6").Value = -78
$ws.Range("S6").Value = 122
$ws.Range("T6").Value = 180
$ws.Range("U6").Value = -222
$ws.Range("V6").Value = 648
$ws.Range("W6").Value = -6.22
$ws.Range("X6").Value = -7.25
$ws.Range("Y6").Value = -22.45
$ws.Range("Z6").Value = -9.98
$ws.Range("AA6").Value = 153.83
$ws.Range("AB6").Value = 48.54
$ws.Range("AC6").Value = -181
$ws.Range("AD6").Value = -14.33
$ws.Range("AE6").Value = 712
$ws.Range("AF6").Value = 3.64
$ws.Range("AI6").Value = 0
$ws.Range("AJ6").Value = 77124820
$ws.Range("AG6").ClearContents()
$ws.Range("AH6").ClearContents()

# Row 7: clear all data columns, keep A/B/C labels
$ws.Range("D7:AJ7").ClearContents()

# Row 8: clear all data columns, keep A/B/C labels
$ws.Range("D8:AJ8").ClearContents()

# Row 9: clear all data columns, keep A/B/C labels
$ws.Range("D9:AJ9").ClearContents()
